$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.097.23'
$ws.Range("E2").Value = '  -2.24%  '

$ws.Range("D3").Value = '3.566.07'
$ws.Range("E3").Value = '  -3.36%  '

$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").Value = '619.71'
$ws.Range("E5").Value = '  -7.19%  '

$ws.Range("D6").Value = '154.19'
$ws.Range("E6").Value = '  -3.88%  '

$ws.Range("D7").Value = '3.560.96'

$ws.Range("E8").Value = '  +0.11%  '

$ws.Range("E9").Value = '  -2.26%  '

$ws.Range("E10").Value = '  -3.17%  '

$ws.Range("D11").Value = '6.96'
$ws.Range("E11").Value = '  -2.28%  '

$ws.Range("E12").Value = '  -1.77%  '

$ws.Range("E13").Value = '  -3.53%  '

$ws.Range("D14").Value = '4.169.18'
$ws.Range("E14").Value = '  -3.30%  '

$ws.Range("D15").Value = '32.09'
$ws.Range("E15").Value = '  -2.46%  '

$ws.Range("D16").Value = '3.573.75'
$ws.Range("E16").Value = '  -2.54%  '

$ws.Range("D17").Value = '68.096.91'
$ws.Range("E17").Value = '  -2.23%  '

$ws.Range("E18").Value = '  -0.85%  '

$ws.Range("D19").Value = '6.44'
$ws.Range("E19").Value = '  -0.44%  '

$ws.Range("D20").Value = '15.65'
$ws.Range("E20").Value = '  -3.24%  '

$ws.Range("D21").Value = '459.46'
$ws.Range("E21").Value = '  -2.31%  '

$ws.Range("E22").Value = '  -0.71%  '

$ws.Range("D23").Value = '0.647'
$ws.Range("E23").Value = '  +0.04%  '

$ws.Range("D24").Value = '77.90'
$ws.Range("E24").Value = '  -2.52%  '

$ws.Range("D25").Value = '3.707.40'
$ws.Range("E25").Value = '  -3.34%  '

$ws.Range("E26").Value = '  +0.09%  '

$ws.Range("D27").Value = '10.69'
$ws.Range("E27").Value = '  -2.30%  '

$ws.Range("E28").Value = '  -8.65%  '

$ws.Range("D29").Value = '8.35'
$ws.Range("E29").Value = '  -7.77%  '

$ws.Range("E30").Value = '  -3.50%  '

$ws.Range("D31").Value = '1.63'
$ws.Range("E31").Value = '  -4.02%  '

$ws.Range("E32").Value = '  -0.03%  '

$ws.Range("D33").Value = '26.05'
$ws.Range("E33").Value = '  -2.67%  '

$ws.Range("B35").Value = 'Kaspa'
$ws.Range("C35").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D35").Value = '0.158'
$ws.Range("E35").Value = '  -4.29%  '

$ws.Range("B36").Value = 'RenzoRestakedETH'
$ws.Range("C36").Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range("D36").Value = '3.564.80'
$ws.Range("E36").Value = '  -3.19%  '

$ws.Range("B37").Value = 'NEARProtocol'
$ws.Range("C37").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D37").Value = '6.18'
$ws.Range("E37").Value = '  -4.77%  '

$ws.Range("D38").Value = '8.11'
$ws.Range("E38").Value = '  -4.25%  '

$ws.Range("E39").Value = '  +0.01%  '

$ws.Range("D40").Value = '178.61'
$ws.Range("E40").Value = '  +0.98%  '

$ws.Range("D41").Value = '0.999'
$ws.Range("E41").Value = '  -0.08%  '

$ws.Range("D42").Value = '0.0887'
$ws.Range("E42").Value = '  -2.11%  '

$ws.Range("D43").Value = '5.63'
$ws.Range("E43").Value = '  -7.98%  '

$ws.Range("E44").Value = '  -6.02%  '

$ws.Range("D45").Value = '0.895'
$ws.Range("E45").Value = '  -4.21%  '

$ws.Range("B46").Value = 'OKB'
$ws.Range("C46").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D46").Value = '46.02'
$ws.Range("E46").Value = '  -2.14%  '

$ws.Range("B47").Value = 'InjectiveProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D47").Value = '28.70'
$ws.Range("E47").Value = '  +3.55%  '

$ws.Range("D48").Value = '2.58'
$ws.Range("E48").Value = '  -6.63%  '

$ws.Range("E49").Value = '  -1.44%  '

$ws.Range("D50").Value = '1.21'
$ws.Range("E50").Value = '  -6.06%  '

$ws.Range("E51").Value = '  -5.44%  '

